$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = 44321
$ws.Cells.Item(2, 10).Value = 38
$ws.Cells.Item(2, 11).Value = 15000
$ws.Cells.Item(2, 12).Value = 15000
$ws.Cells.Item(2, 13).Value = 15000
$ws.Cells.Item(2, 16).Value = 1000
$ws.Cells.Item(3, 4).Value = 44344
$ws.Cells.Item(3, 10).Value = 40
$ws.Cells.Item(3, 11).Value = 20000
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = 20000
$ws.Cells.Item(3, 16).Value = 1333
$ws.Cells.Item(4, 4).Value = 44448
$ws.Cells.Item(4, 10).Value = 85
$ws.Cells.Item(4, 11).Value = 21000
$ws.Cells.Item(4, 12).Value = 22000
$ws.Cells.Item(4, 13).Value = 21529
$ws.Cells.Item(4, 16).Value = 1435
$ws.Cells.Item(5, 4).Value = 44327
$ws.Cells.Item(5, 10).Value = 35
$ws.Cells.Item(5, 11).Value = 15000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 15000
$ws.Cells.Item(5, 16).Value = 1000
$ws.Cells.Item(6, 4).Value = 44309
$ws.Cells.Item(6, 10).Value = 50
$ws.Cells.Item(7, 4).Value = 44308
$ws.Cells.Item(7, 11).Value = 16000
$ws.Cells.Item(7, 12).Value = 16000
$ws.Cells.Item(7, 13).Value = 16000
$ws.Cells.Item(7, 16).Value = 1067
$ws.Cells.Item(8, 4).Value = 44340
$ws.Cells.Item(8, 10).Value = 47
$ws.Cells.Item(8, 11).Value = 14000
$ws.Cells.Item(8, 12).Value = 14000
$ws.Cells.Item(8, 13).Value = 14000
$ws.Cells.Item(8, 16).Value = 933
$ws.Cells.Item(9, 4).Value = 44377
$ws.Cells.Item(9, 10).Value = 80
$ws.Cells.Item(9, 11).Value = 18000
$ws.Cells.Item(9, 12).Value = 19000
$ws.Cells.Item(9, 13).Value = 18500
$ws.Cells.Item(9, 16).Value = 1233
$ws.Cells.Item(10, 4).Value = 44314
$ws.Cells.Item(10, 10).Value = 45
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 15000
$ws.Cells.Item(10, 13).Value = 15000
$ws.Cells.Item(10, 16).Value = 1000
$ws.Cells.Item(11, 4).Value = 44406
$ws.Cells.Item(12, 4).Value = 44341
$ws.Cells.Item(12, 10).Value = 40
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 15000
$ws.Cells.Item(12, 13).Value = 15000
$ws.Cells.Item(12, 16).Value = 1000
$ws.Cells.Item(13, 4).Value = 44323
$ws.Cells.Item(13, 11).Value = 15000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 15000
$ws.Cells.Item(13, 16).Value = 1000
$ws.Cells.Item(14, 4).Value = 44328
$ws.Cells.Item(14, 10).Value = 38
$ws.Cells.Item(15, 4).Value = 44316
$ws.Cells.Item(15, 10).Value = 45
$ws.Cells.Item(15, 11).Value = 14000
$ws.Cells.Item(15, 13).Value = 14444
$ws.Cells.Item(15, 16).Value = 963
$ws.Cells.Item(16, 4).Value = 44322
$ws.Cells.Item(16, 10).Value = 70
$ws.Cells.Item(16, 13).Value = 14500
$ws.Cells.Item(16, 16).Value = 967
$ws.Cells.Item(17, 4).Value = 44452
$ws.Cells.Item(17, 10).Value = 73
$ws.Cells.Item(17, 11).Value = 22000
$ws.Cells.Item(17, 12).Value = 23000
$ws.Cells.Item(17, 13).Value = 22479
$ws.Cells.Item(17, 16).Value = 1499
$ws.Cells.Item(18, 4).Value = 44320
$ws.Cells.Item(18, 10).Value = 40
$ws.Cells.Item(19, 4).Value = 44313
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 14000
$ws.Cells.Item(19, 12).Value = 14000
$ws.Cells.Item(19, 13).Value = 14000
$ws.Cells.Item(19, 16).Value = 933
$ws.Cells.Item(20, 4).Value = 44333
$ws.Cells.Item(20, 10).Value = 35
$ws.Cells.Item(21, 4).Value = 44326
$ws.Cells.Item(21, 10).Value = 45
$ws.Cells.Item(21, 11).Value = 15000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 13).Value = 15000
$ws.Cells.Item(21, 16).Value = 1000
$ws.Cells.Item(22, 4).Value = 44438
$ws.Cells.Item(22, 10).Value = 75
$ws.Cells.Item(22, 11).Value = 19000
$ws.Cells.Item(22, 12).Value = 20000
$ws.Cells.Item(22, 13).Value = 19467
$ws.Cells.Item(22, 16).Value = 1298
$ws.Cells.Item(23, 4).Value = 44411
$ws.Cells.Item(23, 10).Value = 50
$ws.Cells.Item(23, 11).Value = 22000
$ws.Cells.Item(23, 12).Value = 22000
$ws.Cells.Item(23, 13).Value = 22000
$ws.Cells.Item(23, 16).Value = 1467
$ws.Cells.Item(24, 4).Value = 44343
$ws.Cells.Item(24, 11).Value = 15000
$ws.Cells.Item(24, 12).Value = 15000
$ws.Cells.Item(24, 13).Value = 15000
$ws.Cells.Item(24, 16).Value = 1000
$ws.Cells.Item(25, 4).Value = 44315
$ws.Cells.Item(25, 10).Value = 65
$ws.Cells.Item(25, 12).Value = 15000
$ws.Cells.Item(25, 13).Value = 14538
$ws.Cells.Item(25, 16).Value = 969
$ws.Cells.Item(26, 4).Value = 44336
$ws.Cells.Item(26, 10).Value = 65
$ws.Cells.Item(26, 11).Value = 14000
$ws.Cells.Item(26, 13).Value = 14462
$ws.Cells.Item(26, 16).Value = 964
$ws.Cells.Item(27, 4).Value = 44334
$ws.Cells.Item(27, 10).Value = 50
$ws.Cells.Item(27, 11).Value = 14000
$ws.Cells.Item(27, 12).Value = 14000
$ws.Cells.Item(27, 13).Value = 14000
$ws.Cells.Item(27, 16).Value = 933
$ws.Cells.Item(28, 4).Value = 44319
$ws.Cells.Item(28, 10).Value = 50
$ws.Cells.Item(28, 11).Value = 15000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 15000
$ws.Cells.Item(28, 16).Value = 1000
$ws.Cells.Item(29, 4).Value = 44455
$ws.Cells.Item(29, 10).Value = 35
$ws.Cells.Item(29, 11).Value = 22000
$ws.Cells.Item(29, 12).Value = 22000
$ws.Cells.Item(29, 13).Value = 22000
$ws.Cells.Item(29, 16).Value = 1467
$ws.Cells.Item(30, 4).Value = 44397
$ws.Cells.Item(30, 10).Value = 73
$ws.Cells.Item(30, 11).Value = 21000
$ws.Cells.Item(30, 12).Value = 22000
$ws.Cells.Item(30, 13).Value = 21521
$ws.Cells.Item(30, 16).Value = 1435
$ws.Cells.Item(31, 4).Value = 44329
$ws.Cells.Item(31, 10).Value = 35
$ws.Cells.Item(31, 11).Value = 15000
$ws.Cells.Item(31, 12).Value = 15000
$ws.Cells.Item(31, 13).Value = 15000
$ws.Cells.Item(31, 16).Value = 1000
$ws.Cells.Item(32, 4).Value = 44330
$ws.Cells.Item(32, 10).Value = 30
$ws.Cells.Item(32, 11).Value = 15000
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 13).Value = 15000
$ws.Cells.Item(32, 16).Value = 1000
$ws.Cells.Item(33, 4).Value = 44370
$ws.Cells.Item(33, 10).Value = 50
$ws.Cells.Item(33, 11).Value = 18000
$ws.Cells.Item(33, 12).Value = 18000
$ws.Cells.Item(33, 13).Value = 18000
$ws.Cells.Item(33, 16).Value = 1200
$ws.Cells.Item(34, 4).Value = 44312
$ws.Cells.Item(34, 10).Value = 80
$ws.Cells.Item(34, 11).Value = 13000
$ws.Cells.Item(34, 12).Value = 14000
$ws.Cells.Item(34, 13).Value = 13562
$ws.Cells.Item(34, 16).Value = 904
$ws.Cells.Item(35, 4).Value = 44399
$ws.Cells.Item(35, 11).Value = 22000
$ws.Cells.Item(35, 12).Value = 22000
$ws.Cells.Item(35, 13).Value = 22000
$ws.Cells.Item(35, 16).Value = 1467
